$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. Excel shifts every existing column
# from B onward one position to the right (B->C, C->D, ... K->L).
$ws.Columns("B:B").Insert()

# Row 1 is the header row: the new B1 becomes the newest price-check
# timestamp (everything else already shifted right, so the former B1..K1
# values are now sitting in C1..L1).
$ws.Range("B1").Value = "2025-12-19 20:41"

# Rows 52-76 hold the per-SKU rolling price history (columns B..G are the
# 6 most-recent price snapshots). After the column insert, the new B is
# blank and the old "current price" value now lives in C - copy it back
# into B so the newest snapshot again shows the latest known price.
for ($r = 52; $r -le 76; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 3).Value2
}
